$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.680.34"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "3.622.76"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Formula = "'583.72"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").Formula = "'175.68"
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("E7").Value = "  +3.52%  "
$ws.Range("D8").Value = "3.616.29"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Formula = "'0.194"
$ws.Range("E10").Value = "  -4.89%  "
$ws.Range("D11").Formula = "'6.65"
$ws.Range("E11").Value = "  +13.72%  "
$ws.Range("D12").Formula = "'0.619"
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("D13").Formula = "'48.38"
$ws.Range("E13").Value = "  -3.38%  "
$ws.Range("D14").Formula = "'0.0000282"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "4.206.71"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Formula = "'672.06"
$ws.Range("E16").Value = "  -3.33%  "
$ws.Range("D17").Formula = "'9.00"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "3.621.74"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").Value = "70.674.99"
$ws.Range("E19").Value = "  -2.08%  "
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").Formula = "'17.77"
$ws.Range("E21").Value = "  -4.32%  "
$ws.Range("D22").Formula = "'11.45"
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("D23").Formula = "'0.939"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").Formula = "'17.03"
$ws.Range("E24").Value = "  -5.51%  "
$ws.Range("D25").Formula = "'99.74"
$ws.Range("E25").Value = "  -3.99%  "
$ws.Range("D26").Formula = "'3.90"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("E27").Value = "  -3.53%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Formula = "'9.80"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").Formula = "'34.46"
$ws.Range("E30").Value = "  -2.62%  "
$ws.Range("D31").Formula = "'9.15"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Formula = "'3.28"
$ws.Range("E32").Value = "  -4.12%  "
$ws.Range("D33").Formula = "'7.58"
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("D34").Formula = "'1.37"
$ws.Range("E34").Value = "  -5.12%  "
$ws.Range("D35").Formula = "'3.95"
$ws.Range("E35").Value = "  -5.36%  "
$ws.Range("D36").Formula = "'571.78"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("D37").Formula = "'11.07"
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("D39").Formula = "'58.40"
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("D40").Formula = "'0.999"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Formula = "'0.0453"
$ws.Range("E41").Value = "  -2.24%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Formula = "'0.346"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.534.97"
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Formula = "'0.140"
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("D45").Formula = "'34.25"
$ws.Range("E45").Value = "  -4.90%  "
$ws.Range("D46").Value = "0.0₃0728"
$ws.Range("E46").Value = "  -5.92%  "
$ws.Range("D47").Formula = "'2.99"
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("D48").Formula = "'2.66"
$ws.Range("E48").Value = "  -4.41%  "
$ws.Range("D49").Formula = "'0.135"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").Formula = "'137.10"
$ws.Range("E50").Value = "  +3.34%  "
$ws.Range("E51").Value = "  -5.55%  "